$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Article Title:  ___...___" paragraph -> insert the bold
# text "life is good" after the first three underscores, so the single
# run splits into three runs (same bold/size formatting throughout):
#   "Article Title:  ___" | "life is good" | "___....___"
# ---------------------------------------------------------------------
$titleFind = $d.Content.Duplicate
$titleFind.Find.Execute("Article Title:  ___", $true, $false, $false, $false,
                         $false, $true, 1, $false, "", 0) | Out-Null
$insertPos = $titleFind.End

$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter("life is good")

# Re-read the just-inserted span as its own Range and round-trip its
# FormattedText back onto itself. That forces the engine to materialize
# a standalone run boundary at [insertPos, insertPos+Len) instead of
# silently folding the new text back into the (identically formatted)
# neighboring runs.
$newTextRange = $d.Range($insertPos, $insertPos + "life is good".Length)
$newTextRange.FormattedText = $newTextRange.FormattedText

# ---------------------------------------------------------------------
# Change 2: "Question # 2:  ___...___" paragraph currently is three
# runs ("Question # " / "2" / ":  ___...___") -> collapse into a single
# run, matching how "Question # 1:  ___...___" (already one run) looks.
# ---------------------------------------------------------------------
$q1Find = $d.Content.Duplicate
$q1Find.Find.Execute("Question # 1:*____", $true, $false, $true, $false,
                      $false, $true, 1, $false, "", 0) | Out-Null
$q1Para = $q1Find.Paragraphs(1).Range

$q2Find = $d.Content.Duplicate
$q2Find.Find.Execute("Question # 2:*____", $true, $false, $true, $false,
                      $false, $true, 1, $false, "", 0) | Out-Null
$q2Para = $q2Find.Paragraphs(1).Range

# Borrow the fully-formed FormattedText from the Question # 1 run (a
# single run already carrying the exact target rPr), stamp it onto the
# Question # 2 paragraph (collapsing its 3 runs into 1), then fix the
# digit back up with a same-length in-place Find/Replace.
$q2Para.FormattedText = $q1Para.FormattedText

$q2Fix = $d.Content.Duplicate
$q2Fix.Find.Execute("Question # 1:  ____________________________________________________________________________________________________________________________________________________________",
                     $true, $false, $false, $false, $false, $true, 1, $false,
                     "Question # 2:  ____________________________________________________________________________________________________________________________________________________________", 2) | Out-Null
